$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1): B1=B, C1=C, D1=D, E1=E
$ws.Range("D1").Value = "D"
$ws.Range("E1").Value = "E"

# Row labels (column A): A2=B, A3=C, A4=D, A5=E
$ws.Range("A4").Value = "D"
$ws.Range("A5").Value = "E"

# Clear any special formatting on C2 (reset to default/general style)
$ws.Range("C2").ClearFormats()

# Fill in the pairwise comparison matrix values
# Row 2 (B)
$ws.Range("C2").Value = 0.2
$ws.Range("D2").Value = 0.2
$ws.Range("E2").Value = 0.2

# Row 3 (C)
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.2
$ws.Range("E3").Value = 0.2

# Row 4 (D)
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.2

# Row 5 (E)
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 1

# Update the active selection to match the target (D10)
$ws.Range("D10").Select()
